# UPDATE technology portfolios for Norway
$wb = $excel.ActiveWorkbook

# Update the base "2025" sheet with the new input values.
# Other sheets (2030-2050) derive D2/E2 via formulas referencing '2025'!D2 / '2025'!E2,
# but C2 is a static value duplicated on every sheet and must be updated everywhere.
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("C2").Value = 8.5787499999999994
$ws2025.Range("D2").Value = 416394
$ws2025.Range("E2").Value = 360000

foreach ($name in @("2030", "2035", "2040", "2045", "2050")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("C2").Value = 8.5787499999999994
}

$excel.Calculate()
